$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2400.1667
$ws.Range("I40").Value = 1450.5
$ws.Range("J40").Value = 2875
$ws.Range("K40").Value = 1450.5
$ws.Range("L40").Value = 2875
$ws.Range("M40").Value = -1275.5
$ws.Range("N40").Value = -3225
$ws.Range("H58").Value = 10874790
$ws.Range("J58").Value = 16668833
$ws.Range("L58").Value = 50006499
$ws.Range("N58").Value = -50006799
$ws.Range("H61").Value = 4787040
$ws.Range("I61").Value = 8333862
$ws.Range("J61").Value = 57944.668
$ws.Range("K61").Value = 25001586
$ws.Range("L61").Value = 173834.004
$ws.Range("M61").Value = -25001414
$ws.Range("N61").Value = -174178.004
$ws.Range("H64").Value = 93836.09
$ws.Range("J64").Value = 3733.3333
$ws.Range("L64").Value = 3733.3333
$ws.Range("N64").Value = -4229.3333
$ws.Range("H67").Value = 93836.09
$ws.Range("J67").Value = 3733.3333
$ws.Range("L67").Value = 3733.3333
$ws.Range("N67").Value = -5449.3333
$ws.Range("H76").Value = 5575
$ws.Range("J76").Value = 6668
$ws.Range("L76").Value = 6668
$ws.Range("N76").Value = -7298
$ws.Range("H79").Value = 5575
$ws.Range("J79").Value = 6668
$ws.Range("L79").Value = 6668
$ws.Range("N79").Value = -8852
$ws.Range("H82").Value = 3085.2222
$ws.Range("I82").Value = 1262.3334
$ws.Range("J82").Value = 3996.6667
$ws.Range("K82").Value = 3787.0002
$ws.Range("L82").Value = 11990.0001
$ws.Range("M82").Value = -3381.0002
$ws.Range("N82").Value = -12802.0001
$ws.Range("H85").Value = 3085.2222
$ws.Range("I85").Value = 1262.3334
$ws.Range("J85").Value = 3996.6667
$ws.Range("K85").Value = 3787.0002
$ws.Range("L85").Value = 11990.0001
$ws.Range("M85").Value = -2383.0002
$ws.Range("N85").Value = -14798.0001
$ws.Range("H118").Value = 7318
$ws.Range("I118").Value = 10342
$ws.Range("J118").Value = 1270
$ws.Range("K118").Value = 31026
$ws.Range("L118").Value = 3810
$ws.Range("M118").Value = -29369
$ws.Range("N118").Value = -7124
$ws.Range("H129").Value = 2531.1187
$ws.Range("J129").Value = 905.8
$ws.Range("L129").Value = 2717.4
$ws.Range("N129").Value = -12717.4
$ws.Range("H135").Value = 3597.4285
$ws.Range("I135").Value = 829.64703
$ws.Range("K135").Value = 7466.82327
$ws.Range("M135").Value = -4931.82327
$ws.Range("H137").Value = 1626.3235
$ws.Range("I137").Value = 1245.8462
$ws.Range("K137").Value = 3737.5386
$ws.Range("M137").Value = -1187.5386

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31007.871
$ws.Range("I32").Value = 11470.328
$ws.Range("J32").Value = 120322.36
$ws.Range("K32").Value = 11470.328
$ws.Range("L32").Value = 120322.36
$ws.Range("M32").Value = -11183.328
$ws.Range("N32").Value = -120896.36
$ws.Range("H61").Value = 2068
$ws.Range("I61").Value = 1750
$ws.Range("J61").Value = 2110.4
$ws.Range("K61").Value = 1750
$ws.Range("L61").Value = 2110.4
$ws.Range("M61").Value = -1538
$ws.Range("N61").Value = -2534.4
$ws.Range("H132").Value = 18463.514
$ws.Range("I132").Value = 22632.828
$ws.Range("K132").Value = 67898.484
$ws.Range("M132").Value = -65368.484
$ws.Range("H136").Value = 2068
$ws.Range("I136").Value = 1750
$ws.Range("J136").Value = 2110.4
$ws.Range("K136").Value = 5250
$ws.Range("L136").Value = 6331.200000000001
$ws.Range("M136").Value = -2700
$ws.Range("N136").Value = -11431.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 33859
$ws.Range("J101").Value = 33859
$ws.Range("L101").Value = 33859
$ws.Range("N101").Value = -40349
$ws.Range("H134").Value = 3250.2954
$ws.Range("I134").Value = 3130
$ws.Range("J134").Value = 3718.111
$ws.Range("K134").Value = 9390
$ws.Range("L134").Value = 11154.333
$ws.Range("M134").Value = -6855
$ws.Range("N134").Value = -16224.333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2892
$ws.Range("I99").Value = 2860
$ws.Range("J99").Value = 2900
$ws.Range("K99").Value = 2860
$ws.Range("L99").Value = 2900
$ws.Range("M99").Value = -1362
$ws.Range("N99").Value = -5896
$ws.Range("H122").Value = 1006.3333
$ws.Range("J122").Value = 1010
$ws.Range("L122").Value = 3030
$ws.Range("N122").Value = -7930
$ws.Range("H126").Value = 2892
$ws.Range("I126").Value = 2860
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 8580
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -6110
$ws.Range("N126").Value = -13640

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 810.5
$ws.Range("I113").Value = 1023.3684
$ws.Range("J113").Value = 572.58826
$ws.Range("K113").Value = 3070.1052
$ws.Range("L113").Value = 1717.76478
$ws.Range("M113").Value = -900.1052
$ws.Range("N113").Value = -6057.76478
$ws.Range("H131").Value = 789.74
$ws.Range("I131").Value = 455.05
$ws.Range("J131").Value = 873.4125
$ws.Range("K131").Value = 1365.15
$ws.Range("L131").Value = 2620.2375
$ws.Range("M131").Value = 3674.85
$ws.Range("N131").Value = -12700.2375
$ws.Range("H134").Value = 3745.0715
$ws.Range("I134").Value = 2033.579
$ws.Range("J134").Value = 7358.222
$ws.Range("K134").Value = 6100.737
$ws.Range("L134").Value = 22074.666
$ws.Range("M134").Value = -1030.737
$ws.Range("N134").Value = -32214.666

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5885342
$ws.Range("I126").Value = 3159.8
$ws.Range("J126").Value = 11767525
$ws.Range("K126").Value = 9479.400000000001
$ws.Range("L126").Value = 35302575
$ws.Range("M126").Value = -7009.400000000001
$ws.Range("N126").Value = -35307515
$ws.Range("H132").Value = 3301.64
$ws.Range("I132").Value = 2364.125
$ws.Range("J132").Value = 4968.3335
$ws.Range("K132").Value = 7092.375
$ws.Range("L132").Value = 14905.0005
$ws.Range("M132").Value = -4562.375
$ws.Range("N132").Value = -19965.0005
$ws.Range("H141").Value = 49433.332
$ws.Range("J141").Value = 49433.332
$ws.Range("L141").Value = 49433.332
$ws.Range("N141").Value = -59793.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2534.9167
$ws.Range("I7").Value = 1777.3334
$ws.Range("J7").Value = 3292.5
$ws.Range("K7").Value = 1777.3334
$ws.Range("L7").Value = 3292.5
$ws.Range("M7").Value = -1665.3334
$ws.Range("N7").Value = -3516.5
$ws.Range("H40").Value = 69910.92999999999
$ws.Range("I40").Value = 337533
$ws.Range("J40").Value = 3005.4167
$ws.Range("K40").Value = 337533
$ws.Range("L40").Value = 3005.4167
$ws.Range("M40").Value = -337397
$ws.Range("N40").Value = -3277.4167
$ws.Range("H102").Value = 45000
$ws.Range("J102").Value = 45000
$ws.Range("L102").Value = 45000
$ws.Range("N102").Value = -51490
$ws.Range("H126").Value = 2534.9167
$ws.Range("I126").Value = 1777.3334
$ws.Range("J126").Value = 3292.5
$ws.Range("K126").Value = 5332.0002
$ws.Range("L126").Value = 9877.5
$ws.Range("M126").Value = -2862.0002
$ws.Range("N126").Value = -14817.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1971.174
$ws.Range("I126").Value = 1690.7059
$ws.Range("J126").Value = 2765.8333
$ws.Range("K126").Value = 5072.1177
$ws.Range("L126").Value = 8297.499899999999
$ws.Range("M126").Value = -2602.1177
$ws.Range("N126").Value = -13237.4999
$ws.Range("H140").Value = 55866.668
$ws.Range("J140").Value = 55866.668
$ws.Range("L140").Value = 55866.668
$ws.Range("N140").Value = -66226.66800000001
$ws.Range("H141").Value = 57537.855
$ws.Range("J141").Value = 57537.855
$ws.Range("L141").Value = 57537.855
$ws.Range("N141").Value = -67897.85500000001

Write-Host "Applied 200 cell updates across 8 sheets"
